# Trading update: 2026-02-17 13:10:15
# Appends the newest (still-open) MarketMaking trade as row 15 to both the
# "All Trades" log and the per-strategy "MarketMaking" log.

$wb = $excel.ActiveWorkbook

function Add-TradeRow15($sheetName) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A15").Value = 14

    # Column B holds a literal "YYYY-MM-DD" *text* label (not a real Excel
    # date) in every existing row, so force Text formatting before writing
    # it - otherwise Excel's automatic date recognition would turn the
    # string into a date serial number.
    $ws.Range("B15").NumberFormat = "@"
    $ws.Range("B15").Value = "2026-02-17"

    $ws.Range("C15").Value = "13:10:15"
    $ws.Range("D15").Value = "MarketMaking"
    $ws.Range("E15").Value = "UP"
    $ws.Range("F15").Value = 0.53
    $ws.Range("G15").Value = ""
    $ws.Range("H15").Value = "OPEN"
    $ws.Range("I15").Value = 0
    $ws.Range("J15").Value = 0
    $ws.Range("K15").Value = 99.6022341356021
    $ws.Range("L15").Value = 0
    $ws.Range("M15").Value = 0
    $ws.Range("N15").Value = 0.6
    $ws.Range("O15").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P15").Value = ""
    $ws.Range("Q15").Value = 0

    # Reset to the workbook's default style so the new row doesn't pick up
    # a stray number-format style (only the literal text content matters).
    $ws.Range("A15:Q15").Style = "Normal"
}

Add-TradeRow15 "All Trades"
Add-TradeRow15 "MarketMaking"
